$d = $word.ActiveDocument

# Locate the paragraph that ends the "Reflection" bullet group:
# "Should intercept methods like Class.getFields() to interface with mirrors instead"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Should intercept methods like Class.getFields*") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not find anchor paragraph"
}

# Insert a new paragraph right after it for "Enums" (second-level bullet, ilvl=1)
$r = $target.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$enumsPara = $target.Next()
$enumsPara.Range.Text = "Enums"
$enumsPara.Range.ListFormat.ListLevelNumber = 2

# Insert another paragraph after "Enums" for the follow-up note (third-level bullet, ilvl=2)
$r2 = $enumsPara.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$flagPara = $enumsPara.Next()
$flagPara.Range.Text = "Just remove the flag – no need?"
$flagPara.Range.ListFormat.ListLevelNumber = 3

Write-Host "Inserted paragraphs. New count:" $d.Paragraphs.Count
